$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1926605504587156
$ws.Range("C2").Value = 0.5626911314984709
$ws.Range("J2").Value = 0.006116207951070336
$ws.Range("P2").Value = 0.1437308868501529
$ws.Range("S2").Value = 0.09480122324159021
$ws.Range("B3").Value = 0.005208333333333333
$ws.Range("C3").Value = 0.02083333333333333
$ws.Range("J3").Value = 0.03125
$ws.Range("P3").Value = 0.7291666666666666
$ws.Range("S3").Value = 0.2135416666666667
$ws.Range("J4").Value = 0.1176470588235294
$ws.Range("P4").Value = 0.5588235294117647
$ws.Range("S4").Value = 0.3235294117647059
$ws.Range("B6").Value = 0.05098039215686274
$ws.Range("D6").Value = 0.007843137254901961
$ws.Range("F6").Value = 0.07843137254901961
$ws.Range("J6").Value = 0.2627450980392157
$ws.Range("O6").Value = 0.01568627450980392
$ws.Range("Q6").Value = 0.1411764705882353
$ws.Range("R6").Value = 0.09411764705882353
$ws.Range("S6").Value = 0.3490196078431372
$ws.Range("B7").Value = 0.0796812749003984
$ws.Range("D7").Value = 0.01593625498007968
$ws.Range("E7").Value = 0.00398406374501992
$ws.Range("F7").Value = 0.05179282868525897
$ws.Range("J7").Value = 0.1354581673306773
$ws.Range("O7").Value = 0.0199203187250996
$ws.Range("Q7").Value = 0.1553784860557769
$ws.Range("R7").Value = 0.1035856573705179
$ws.Range("S7").Value = 0.4342629482071713
$ws.Range("B8").Value = 0.08173913043478261
$ws.Range("D8").Value = 0.008695652173913044
$ws.Range("F8").Value = 0.05391304347826087
$ws.Range("J8").Value = 0.1060869565217391
$ws.Range("O8").Value = 0.02434782608695652
$ws.Range("Q8").Value = 0.1617391304347826
$ws.Range("R8").Value = 0.1304347826086956
$ws.Range("S8").Value = 0.4330434782608696
$ws.Range("B9").Value = 0.06944444444444445
$ws.Range("D9").Value = 0.02777777777777778
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.1180555555555556
$ws.Range("O9").Value = 0.03472222222222222
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.09722222222222222
$ws.Range("S9").Value = 0.4305555555555556
$ws.Range("B10").Value = 0.1197387518142235
$ws.Range("D10").Value = 0.01669085631349782
$ws.Range("E10").Value = 0.001451378809869376
$ws.Range("F10").Value = 0.06894049346879536
$ws.Range("J10").Value = 0.1190130624092888
$ws.Range("O10").Value = 0.02031930333817126
$ws.Range("Q10").Value = 0.1777939042089985
$ws.Range("R10").Value = 0.08925979680696662
$ws.Range("S10").Value = 0.3867924528301887
$ws.Range("G11").Value = 0.1555555555555556
$ws.Range("J11").Value = 0.08641975308641975
$ws.Range("K11").Value = 0.1925925925925926
$ws.Range("L11").Value = 0.528395061728395
$ws.Range("S11").Value = 0.03703703703703703
$ws.Range("G12").Value = 0.7162162162162162
$ws.Range("J12").Value = 0.1801801801801802
$ws.Range("K12").Value = 0.02252252252252252
$ws.Range("L12").Value = 0.03153153153153153
$ws.Range("S12").Value = 0.04954954954954955
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.0625
$ws.Range("G14").Value = 0.8571428571428571
$ws.Range("J14").Value = 0.1428571428571428
$ws.Range("F15").Value = 0.02962962962962963
$ws.Range("H15").Value = 0.1518518518518518
$ws.Range("I15").Value = 0.04444444444444445
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.0962962962962963
$ws.Range("M15").Value = 0.02222222222222222
$ws.Range("N15").Value = 0.003703703703703704
$ws.Range("O15").Value = 0.1
$ws.Range("S15").Value = 0.2185185185185185
$ws.Range("F16").Value = 0.03517587939698492
$ws.Range("H16").Value = 0.135678391959799
$ws.Range("I16").Value = 0.05025125628140704
$ws.Range("J16").Value = 0.407035175879397
$ws.Range("K16").Value = 0.1507537688442211
$ws.Range("M16").Value = 0.01005025125628141
$ws.Range("N16").Value = 0.01005025125628141
$ws.Range("O16").Value = 0.02010050251256281
$ws.Range("S16").Value = 0.1809045226130653
$ws.Range("F17").Value = 0.01839080459770115
$ws.Range("H17").Value = 0.2275862068965517
$ws.Range("I17").Value = 0.05747126436781609
$ws.Range("J17").Value = 0.3954022988505747
$ws.Range("K17").Value = 0.09655172413793103
$ws.Range("M17").Value = 0.01609195402298851
$ws.Range("N17").Value = 0.002298850574712644
$ws.Range("O17").Value = 0.06206896551724138
$ws.Range("S17").Value = 0.1241379310344828
$ws.Range("F18").Value = 0.01153846153846154
$ws.Range("H18").Value = 0.2076923076923077
$ws.Range("I18").Value = 0.0576923076923077
$ws.Range("J18").Value = 0.4384615384615385
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.02307692307692308
$ws.Range("N18").Value = 0.003846153846153846
$ws.Range("O18").Value = 0.04615384615384616
$ws.Range("S18").Value = 0.1115384615384615
$ws.Range("F19").Value = 0.02194148936170213
$ws.Range("H19").Value = 0.238031914893617
$ws.Range("I19").Value = 0.05518617021276596
$ws.Range("J19").Value = 0.3311170212765958
$ws.Range("K19").Value = 0.1269946808510638
$ws.Range("M19").Value = 0.02127659574468085
$ws.Range("N19").Value = 0.001329787234042553
$ws.Range("O19").Value = 0.07779255319148937
$ws.Range("S19").Value = 0.1263297872340426
